$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix hunter perk icons: these cells incorrectly showed "Check" and should show "OK"
$ws.Range("D31").Value = "OK"
$ws.Range("D34").Value = "OK"
$ws.Range("D35").Value = "OK"
$ws.Range("D36").Value = "OK"
$ws.Range("D41").Value = "OK"

# Move selection to D43 (matches the author's cursor position when saving)
$ws.Range("D43").Select()
